$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 17, shifting all existing
# data rows (17-97) down to (19-99). This matches the dimension change
# from A1:R97 to A1:R99 seen in the diff.
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(18).Insert()

# Populate new row 17 with the new weekly entry (Primera quality)
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 44970
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 100112037
$ws.Range("G17").Value = "Cebollín"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 150
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 900
$ws.Range("M17").Value = 900
$ws.Range("N17").Value = "$/paquete 6 unidades"
$ws.Range("O17").Value = "Provincia de Diguillín"
$ws.Range("P17").Value = 150
$ws.Range("Q17").Value = 6
$ws.Range("R17").Value = "Hortaliza"

# Populate new row 18 with the new weekly entry (Segunda quality)
$ws.Range("A18").Value = 7
$ws.Range("B18").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C18").Value = "Ñuble"
$ws.Range("D18").Value = 44970
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 100112037
$ws.Range("G18").Value = "Cebollín"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 150
$ws.Range("K18").Value = 700
$ws.Range("L18").Value = 700
$ws.Range("M18").Value = 700
$ws.Range("N18").Value = "$/paquete 6 unidades"
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 117
$ws.Range("Q18").Value = 6
$ws.Range("R18").Value = "Hortaliza"
